$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.660.68'
$ws.Range('E2').Value = '  +0.29%  '
$ws.Range('D3').Value = '3.142.39'
$ws.Range('E3').Value = '  +6.15%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '201.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.42%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '624.86'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.73%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.04%  '
$ws.Range('B8').Value = 'Dogecoin'
$ws.Range('C8').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.217'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +6.62%  '
$ws.Range('B9').Value = 'XRP'
$ws.Range('C9').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.557'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.468'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.162'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.39%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.39%  '
$ws.Range('D13').Value = '3.716.49'
$ws.Range('E13').Value = '  +6.13%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '29.75'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.97%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000202'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +5.80%  '
$ws.Range('D16').Value = '76.631.55'
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '3.128.39'
$ws.Range('E17').Value = '  +5.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.76'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '9.28'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +6.10%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.76'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +20.37%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '398.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.34%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.57'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.45%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.54'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.19%  '
$ws.Range('B24').Value = 'WrappedeETH'
$ws.Range('C24').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D24').Value = '3.290.16'
$ws.Range('E24').Value = '  +6.37%  '
$ws.Range('B25').Value = 'NEARProtocol'
$ws.Range('C25').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.63'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +7.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '73.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +7.19%  '
$ws.Range('E28').Value = '  +0.09%  '
$ws.Range('E29').Value = '  +3.52%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.997'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.50'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.24%  '
$ws.Range('E32').Value = '  +5.56%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '514.63'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.54%  '
$ws.Range('E34').Value = '  +6.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.132'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +18.40%  '
$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '22.00'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.03%  '
$ws.Range('B37').Value = 'FirstDigitalUSD'
$ws.Range('C37').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.00'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.01%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '163.56'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.03%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '197.41'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.08%  '
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.384'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.28%  '
$ws.Range('E42').Value = '  -6.09%  '
$ws.Range('E43').Value = '  +8.74%  '
$ws.Range('E44').Value = '  +0.04%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.811'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +21.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.29'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +7.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.71'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.11%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '41.73'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.56%  '
$ws.Range('B49').Value = 'dogwifhat'
$ws.Range('C49').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.50'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +7.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.613'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +3.53%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.00'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.80%  '
